$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.183.58'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '1.681.89'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.519'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.55'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.98%  '
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = '1.918.28'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").Value = '1.684.50'
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("E14").Value = '  +1.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.539'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = '27.173.62'
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  +1.46%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  +2.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.94%  '
$ws.Range("E28").Value = '  +0.72%  '
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("D32").Value = '1.572.97'
$ws.Range("E32").Value = '  +5.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("E34").Value = '  +2.41%  '
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.602'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.95%  '
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.936'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.03%  '
$ws.Range("E39").Value = '  -0.38%  '
$ws.Range("E40").Value = '  +4.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.47%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  -4.26%  '
$ws.Range("E44").Value = '  -2.90%  '
$ws.Range("D45").Value = '1.827.98'
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.60%  '
$ws.Range("D49").Value = '0.0₆0108'
$ws.Range("E49").Value = '  +2.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.45%  '
$ws.Range("E51").Value = '  +1.56%  '
